$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 86
$ws.Range("H86").Value = 3677.889
$ws.Range("I86").Value = 3614
$ws.Range("K86").Value = 3614
$ws.Range("M86").Value = -2491
# Row 89
$ws.Range("H89").Value = 3677.889
$ws.Range("I89").Value = 3614
$ws.Range("K89").Value = 18070
$ws.Range("M89").Value = -12454
# Row 97
$ws.Range("H97").Value = 1119.2
$ws.Range("J97").Value = 1119.2
$ws.Range("L97").Value = 3357.6
$ws.Range("N97").Value = -4349.6
# Row 101
$ws.Range("H101").Value = 83335310
$ws.Range("I101").Value = 125000340
$ws.Range("J101").Value = 5250
$ws.Range("K101").Value = 375001020
$ws.Range("L101").Value = 15750
$ws.Range("M101").Value = -374999398
$ws.Range("N101").Value = -18994
# Row 125
$ws.Range("H125").Value = 6538410.5
$ws.Range("I125").Value = 1461
$ws.Range("J125").Value = 7754587
$ws.Range("K125").Value = 13149
$ws.Range("L125").Value = 69791283
$ws.Range("M125").Value = -10689
$ws.Range("N125").Value = -69796203
# Row 137
$ws.Range("H137").Value = 76931.71000000001
$ws.Range("I137").Value = 225159.62
$ws.Range("K137").Value = 675478.86
$ws.Range("M137").Value = -672928.86
# Row 138
$ws.Range("H138").Value = 3300
$ws.Range("I138").Value = 2577.7058
$ws.Range("K138").Value = 7733.117400000001
$ws.Range("M138").Value = -2593.117400000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 8
$ws.Range("H8").Value = 5000
$ws.Range("J8").Value = 5000
$ws.Range("L8").Value = 5000
$ws.Range("N8").Value = -5288
# Row 45
$ws.Range("H45").Value = 5997035
$ws.Range("I45").Value = 11067381
$ws.Range("J45").Value = 4808
$ws.Range("K45").Value = 11067381
$ws.Range("L45").Value = 4808
$ws.Range("M45").Value = -11067004
$ws.Range("N45").Value = -5562
# Row 74
$ws.Range("H74").Value = 75756.13
$ws.Range("I74").Value = 5432.3794
$ws.Range("K74").Value = 5432.3794
$ws.Range("M74").Value = -4558.3794
# Row 77
$ws.Range("H77").Value = 75756.13
$ws.Range("I77").Value = 5432.3794
$ws.Range("K77").Value = 27161.897
$ws.Range("M77").Value = -22793.897

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 4167922.8
$ws.Range("I105").Value = 4465503
$ws.Range("J105").Value = 1800
$ws.Range("K105").Value = 4465503
$ws.Range("L105").Value = 1800
$ws.Range("M105").Value = -4463756
$ws.Range("N105").Value = -5294
# Row 134
$ws.Range("H134").Value = 2761.4043
$ws.Range("I134").Value = 1110.9722
$ws.Range("K134").Value = 3332.9166
$ws.Range("M134").Value = -797.9165999999996

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 13
$ws.Range("H13").Value = 8000
$ws.Range("J13").Value = 8000
$ws.Range("L13").Value = 8000
$ws.Range("N13").Value = -8278
# Row 31
$ws.Range("H31").Value = 24882.77
$ws.Range("I31").Value = 1130.08
$ws.Range("K31").Value = 1130.08
$ws.Range("M31").Value = -835.0799999999999
# Row 34
$ws.Range("H34").Value = 24882.77
$ws.Range("I34").Value = 1130.08
$ws.Range("K34").Value = 1130.08
$ws.Range("M34").Value = -928.0799999999999
# Row 99
$ws.Range("H99").Value = 4249
$ws.Range("I99").Value = 4000
$ws.Range("J99").Value = 4355.7144
$ws.Range("K99").Value = 4000
$ws.Range("L99").Value = 4355.7144
$ws.Range("M99").Value = -2502
$ws.Range("N99").Value = -7351.7144
# Row 122
$ws.Range("H122").Value = 2925.2
$ws.Range("J122").Value = 3221.6
$ws.Range("L122").Value = 9664.799999999999
$ws.Range("N122").Value = -14564.8
# Row 126
$ws.Range("H126").Value = 4249
$ws.Range("I126").Value = 4000
$ws.Range("J126").Value = 4355.7144
$ws.Range("K126").Value = 12000
$ws.Range("L126").Value = 13067.1432
$ws.Range("M126").Value = -9530
$ws.Range("N126").Value = -18007.1432
# Row 134
$ws.Range("H134").Value = 2973.186
$ws.Range("I134").Value = 2407.7878
$ws.Range("K134").Value = 7223.3634
$ws.Range("M134").Value = -4688.3634

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 56
$ws.Range("H56").Value = 10422533
$ws.Range("I56").Value = 10422533
$ws.Range("K56").Value = 10422533
$ws.Range("M56").Value = -10422003
# Row 68
$ws.Range("H68").Value = 1150
$ws.Range("I68").Value = 700
$ws.Range("J68").Value = 1300
$ws.Range("K68").Value = 2100
$ws.Range("L68").Value = 3900
$ws.Range("M68").Value = -1289
$ws.Range("N68").Value = -5522
# Row 71
$ws.Range("H71").Value = 1150
$ws.Range("I71").Value = 700
$ws.Range("J71").Value = 1300
$ws.Range("K71").Value = 6300
$ws.Range("L71").Value = 11700
$ws.Range("M71").Value = -2244
$ws.Range("N71").Value = -19812
# Row 133
$ws.Range("H133").Value = 4831.6665
$ws.Range("I133").Value = 4831.6665
$ws.Range("K133").Value = 14494.9995
$ws.Range("M133").Value = -9434.999500000002
# Row 134
$ws.Range("H134").Value = 1520.762
$ws.Range("I134").Value = 1520.762
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 4562.286
$ws.Range("L134").Value = 0
$ws.Range("N134").Value = 507.7139999999999
$ws.Range("M134").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 41
$ws.Range("H41").Value = 12030
$ws.Range("I41").Value = 5000
$ws.Range("K41").Value = 5000
$ws.Range("M41").Value = -4645
# Row 70
$ws.Range("H70").Value = 33343416
$ws.Range("I70").Value = 33343416
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 33343416
$ws.Range("L70").Value = 0
$ws.Range("N70").Value = -33343146
$ws.Range("M70").ClearContents()
# Row 73
$ws.Range("H73").Value = 33343416
$ws.Range("I73").Value = 33343416
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 33343416
$ws.Range("L73").Value = 0
$ws.Range("N73").Value = -33342480
$ws.Range("M73").ClearContents()
# Row 132
$ws.Range("H132").Value = 2947.2432
$ws.Range("I132").Value = 2485.9678
$ws.Range("J132").Value = 5330.5
$ws.Range("K132").Value = 7457.903399999999
$ws.Range("L132").Value = 15991.5
$ws.Range("M132").Value = -4927.903399999999
$ws.Range("N132").Value = -21051.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 5455911
$ws.Range("J2").Value = 15000
$ws.Range("L2").Value = 15000
$ws.Range("N2").Value = -15224
# Row 38
$ws.Range("H38").Value = 1000000000
$ws.Range("I38").Value = 1000000000
$ws.Range("K38").Value = 1000000000
$ws.Range("M38").Value = -999999590
# Row 82
$ws.Range("H82").Value = 3705370.5
$ws.Range("I82").Value = 6174662
$ws.Range("J82").Value = 1432.6666
$ws.Range("K82").Value = 6174662
$ws.Range("L82").Value = 1432.6666
$ws.Range("M82").Value = -6174301
$ws.Range("N82").Value = -2154.6666
# Row 85
$ws.Range("H85").Value = 3705370.5
$ws.Range("I85").Value = 6174662
$ws.Range("J85").Value = 1432.6666
$ws.Range("K85").Value = 6174662
$ws.Range("L85").Value = 1432.6666
$ws.Range("M85").Value = -6173414
$ws.Range("N85").Value = -3928.6666
# Row 100
$ws.Range("H100").Value = 3150.6206
$ws.Range("I100").Value = 2938.9583
$ws.Range("K100").Value = 2938.9583
$ws.Range("M100").Value = -2397.9583
# Row 106
$ws.Range("H106").Value = 9999.5
$ws.Range("J106").Value = 9999.5
$ws.Range("L106").Value = 9999.5
$ws.Range("N106").Value = -12523.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 27779588
$ws.Range("I107").Value = 33335320
$ws.Range("J107").Value = 930
$ws.Range("K107").Value = 100005960
$ws.Range("L107").Value = 2790
$ws.Range("M107").Value = -100004040
$ws.Range("N107").Value = -6630
# Row 122
$ws.Range("H122").Value = 2183.65
$ws.Range("I122").Value = 1431.1
$ws.Range("J122").Value = 2936.2
$ws.Range("K122").Value = 4293.299999999999
$ws.Range("L122").Value = 8808.599999999999
$ws.Range("M122").Value = -1843.299999999999
$ws.Range("N122").Value = -13708.6
# Row 133
$ws.Range("H133").Value = 75000
$ws.Range("J133").Value = 75000
$ws.Range("L133").Value = 75000
$ws.Range("N133").Value = -85120
# Row 136
$ws.Range("H136").Value = 2241.9375
$ws.Range("I136").Value = 1741.5
$ws.Range("K136").Value = 5224.5
$ws.Range("M136").Value = -2674.5
